$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.417.44'
$ws.Range('E2').Value = '  -1.15%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.502.40'
$ws.Range('E3').Value = '  -2.54%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '606.51'
$ws.Range('E5').Value = '  -3.14%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '149.48'
$ws.Range('E6').Value = '  -4.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.501.72'
$ws.Range('E7').Value = '  -2.51%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  -1.70%  '
$ws.Range('E10').Value = '  -2.71%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.04'
$ws.Range('E11').Value = '  +1.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.424'
$ws.Range('E12').Value = '  -2.46%  '
$ws.Range('E13').Value = '  -3.43%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.092.57'
$ws.Range('E14').Value = '  -2.84%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '31.52'
$ws.Range('E15').Value = '  -1.76%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.498.04'
$ws.Range('E16').Value = '  -3.37%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '67.420.55'
$ws.Range('E17').Value = '  -1.27%  '
$ws.Range('E18').Value = '  -0.44%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.42'
$ws.Range('E19').Value = '  -0.47%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.15'
$ws.Range('E20').Value = '  -3.36%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '445.68'
$ws.Range('E21').Value = '  -2.78%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.14'
$ws.Range('E22').Value = '  -7.00%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.623'
$ws.Range('E23').Value = '  -3.02%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '77.33'
$ws.Range('E24').Value = '  -1.07%  '
$ws.Range('B25').Value = 'WrappedeETH'
$ws.Range('C25').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.640.64'
$ws.Range('E25').Value = '  -2.86%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  +7.60%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.34'
$ws.Range('E28').Value = '  -3.81%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.30'
$ws.Range('E29').Value = '  -1.81%  '
$ws.Range('E30').Value = '  -4.57%  '
$ws.Range('E31').Value = '  +0.00%  '
$ws.Range('E32').Value = '  -7.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.165'
$ws.Range('E33').Value = '  +3.27%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.66'
$ws.Range('E34').Value = '  -1.85%  '
$ws.Range('E35').Value = '  -2.24%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.492.44'
$ws.Range('E36').Value = '  -2.99%  '
$ws.Range('E37').Value = '  -4.87%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '8.05'
$ws.Range('E38').Value = '  -1.59%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  -0.16%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '177.45'
$ws.Range('E41').Value = '  +0.21%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.18'
$ws.Range('E42').Value = '  +1.44%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0871'
$ws.Range('E43').Value = '  -1.07%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.41'
$ws.Range('E44').Value = '  -4.12%  '
$ws.Range('E45').Value = '  -2.75%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '45.35'
$ws.Range('E46').Value = '  -1.70%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '27.32'
$ws.Range('E47').Value = '  -4.83%  '
$ws.Range('E48').Value = '  +2.68%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.54'
$ws.Range('E49').Value = '  -2.75%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.55'
$ws.Range('E50').Value = '  -2.33%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.986'
$ws.Range('E51').Value = '  -2.91%  '
